$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Proyect" column (column B), shifting subsequent columns left
$ws.Columns("B").Delete()

# Fill in the new sample data row (row 2)
$ws.Range("A2").Value = "Miami"
$ws.Range("B2").Value = "Chile"
$ws.Range("C2").Value = "CCU"
$ws.Range("D2").Value = "Fiber"
$ws.Range("E2").Value = "j123"
$ws.Range("F2").Value = "Nodo"
$ws.Range("G2").Value = "un"
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = 30

# Set the active selection to I3
$ws.Range("I3").Select()
